# Insert a new data row at row 224 (pushing the existing rows 224..276
# down to 225..277) and populate the new row with the latest price
# report for "Zapallo italiano" at Vega Monumental Concepción.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 224:276 down one row, just like Excel's own
# right-click > Insert on a row header.
$ws.Rows.Item(224).Insert()

# Fill in the newly inserted (currently blank) row 224.
$ws.Range("A224").Value = 11
$ws.Range("B224").Value = "Vega Monumental Concepción"
$ws.Range("C224").Value = "Bíobío"
$ws.Range("D224").Value = 45204
$ws.Range("E224").Value = 8
$ws.Range("F224").Value = 100112032
$ws.Range("G224").Value = "Zapallo italiano"
$ws.Range("H224").Value = "Sin especificar"
$ws.Range("I224").Value = "Primera"
$ws.Range("J224").Value = 100
$ws.Range("K224").Value = 19000
$ws.Range("L224").Value = 20000
$ws.Range("M224").Value = 19500
$ws.Range("N224").Value = '$/caja 50 unidades'
$ws.Range("O224").Value = "Región de Arica y Parinacota"
$ws.Range("P224").Value = 390
$ws.Range("Q224").Value = 50
$ws.Range("R224").Value = "Hortaliza"
